# Update cryptos list values (prices and 1h-volume deltas) to match
# the latest scrape, mirroring the GitHub Actions commit on
# Wed Jul 26 17:13:33 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.323.19"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.860.56"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'0.7022"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").Value = "'237.85"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D8").Value = "'0.07837"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").Value = "'0.3050"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "'24.67"
$ws.Range("E10").Value = "  +6.49%  "
$ws.Range("D11").Value = "'0.08155"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "1.857.42"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "'5.213"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").Value = "'0.7146"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "'89.20"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "29.336.78"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "'5.793"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'241.42"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007770"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.101.72"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'7.507"
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("D25").Value = "'162.68"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "'8.896"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("D27").Value = "'0.1426"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").Value = "'18.09"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("E29").Value = "  -4.71%  "
$ws.Range("D30").Value = "'1.373"
$ws.Range("E30").Value = "  -4.33%  "
$ws.Range("D31").Value = "'1.474"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").Value = "'4.298"
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("D33").Value = "'4.032"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").Value = "'0.05164"
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("D35").Value = "'1.180"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("D36").Value = "'0.7048"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "'0.9980"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'2.674"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("D39").Value = "'0.01842"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Value = "'2.690"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("D41").Value = "1.171.85"
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("D42").Value = "'0.9178"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").Value = "'6.011"
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("D44").Value = "'71.40"
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("D45").Value = "'0.4240"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "'101.64"
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("D48").Value = "'0.5354"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").Value = "'1.748"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("D50").Value = "'9.147"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").Value = "'6.938"
$ws.Range("E51").Value = "  +0.06%  "
